$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "CSC335"
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = "CSC210"

$ws.Range("A11").Value = "CSC372"
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = "CSC210"

$ws.Range("A12").Value = "CSC473"
$ws.Range("B12").Value = 3
$ws.Range("C12").Value = "CSC345"

$ws.Range("A13").Value = "CSC452"
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = "CSC345, CSC252, CSC352"

$ws.Range("A14").Value = "CHEM"

$ws.Range("A14").Select()
